$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value without Excel re-typing it as a
# number/date and without leaving a stray quote-prefix style behind.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '30.151.12'
$ws.Cells.Item(2, 5).Value = '  -0.55%  '

Set-TextCell 3 4 '1.860.21'
$ws.Cells.Item(3, 5).Value = '  -0.56%  '

$ws.Cells.Item(4, 5).Value = '  -0.01%  '

Set-TextCell 5 4 '233.84'
$ws.Cells.Item(5, 5).Value = '  -0.87%  '

$ws.Cells.Item(6, 5).Value = '  -0.03%  '

Set-TextCell 7 4 '0.4679'
$ws.Cells.Item(7, 5).Value = '  -0.58%  '

Set-TextCell 8 4 '42.79'
$ws.Cells.Item(8, 5).Value = '  -0.47%  '

Set-TextCell 9 4 '0.2847'
$ws.Cells.Item(9, 5).Value = '  -1.33%  '

Set-TextCell 10 4 '0.06467'
$ws.Cells.Item(10, 5).Value = '  -2.32%  '

$ws.Cells.Item(11, 5).Value = '  -3.67%  '

Set-TextCell 12 4 '0.07754'
$ws.Cells.Item(12, 5).Value = '  -3.40%  '

Set-TextCell 13 4 '1.866.14'
$ws.Cells.Item(13, 5).Value = '  -0.26%  '

Set-TextCell 14 4 '93.43'
$ws.Cells.Item(14, 5).Value = '  -4.17%  '

Set-TextCell 15 4 '0.6795'
$ws.Cells.Item(15, 5).Value = '  -1.06%  '

$ws.Cells.Item(16, 5).Value = '  -2.17%  '

Set-TextCell 17 4 '266.38'
$ws.Cells.Item(17, 5).Value = '  -1.82%  '

Set-TextCell 18 4 '30.131.57'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '

$ws.Cells.Item(19, 5).Value = '  -5.90%  '

Set-TextCell 20 4 '0.000007582'
$ws.Cells.Item(20, 5).Value = '  -1.87%  '

$ws.Cells.Item(21, 5).Value = '  +0.02%  '

Set-TextCell 22 4 '2.116.80'
$ws.Cells.Item(22, 5).Value = '  +0.08%  '

$ws.Cells.Item(23, 5).Value = '  +0.03%  '

Set-TextCell 24 4 '5.125'
$ws.Cells.Item(24, 5).Value = '  -3.43%  '

Set-TextCell 25 4 '6.093'
$ws.Cells.Item(25, 5).Value = '  -2.11%  '

$ws.Cells.Item(26, 5).Value = '  -1.05%  '

Set-TextCell 27 4 '165.08'
$ws.Cells.Item(27, 5).Value = '  -2.18%  '

Set-TextCell 28 4 '18.48'
$ws.Cells.Item(28, 5).Value = '  -2.52%  '

Set-TextCell 29 4 '1.879'
$ws.Cells.Item(29, 5).Value = '  -4.04%  '

$ws.Cells.Item(30, 5).Value = '  -0.92%  '

Set-TextCell 31 4 '0.09913'
$ws.Cells.Item(31, 5).Value = '  +0.11%  '

$ws.Cells.Item(32, 5).Value = '  -1.42%  '

Set-TextCell 33 4 '4.203'
$ws.Cells.Item(33, 5).Value = '  -4.00%  '

Set-TextCell 34 4 '3.980'
$ws.Cells.Item(34, 5).Value = '  -2.53%  '

Set-TextCell 35 4 '0.04652'
$ws.Cells.Item(35, 5).Value = '  -1.27%  '

$ws.Cells.Item(36, 5).Value = '  -2.04%  '

Set-TextCell 37 4 '0.6866'
$ws.Cells.Item(37, 5).Value = '  -2.29%  '

Set-TextCell 38 4 '2.713'
$ws.Cells.Item(38, 5).Value = '  +0.34%  '

Set-TextCell 39 4 '0.01829'
$ws.Cells.Item(39, 5).Value = '  -2.93%  '

Set-TextCell 40 4 '2.748'
$ws.Cells.Item(40, 5).Value = '  +3.75%  '

Set-TextCell 41 4 '6.280'
$ws.Cells.Item(41, 5).Value = '  -0.44%  '

Set-TextCell 42 4 '71.08'
$ws.Cells.Item(42, 5).Value = '  -2.50%  '

Set-TextCell 43 4 '1.000'
$ws.Cells.Item(43, 5).Value = '  -0.05%  '

Set-TextCell 44 4 '0.8305'
$ws.Cells.Item(44, 5).Value = '  -1.57%  '

Set-TextCell 45 4 '1.880'
$ws.Cells.Item(45, 5).Value = '  -4.11%  '

Set-TextCell 46 4 '101.91'
$ws.Cells.Item(46, 5).Value = '  -1.33%  '

Set-TextCell 47 4 '0.4038'
$ws.Cells.Item(47, 5).Value = '  -3.31%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 48 4 '9.145'
$ws.Cells.Item(48, 5).Value = '  -0.67%  '

$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 49 4 '925.24'
$ws.Cells.Item(49, 5).Value = '  -0.19%  '

Set-TextCell 50 4 '6.924'
$ws.Cells.Item(50, 5).Value = '  -2.32%  '

Set-TextCell 51 4 '33.95'
$ws.Cells.Item(51, 5).Value = '  -1.57%  '
